{"js": "// Apply the \"Added many more features\" edit to the Jackpot Raiders review.\n// Each change is a straightforward whole-paragraph text replacement;\n// we locate the exact original text with a body-wide search (matching\n// the full run text so we only ever hit the intended run) and overwrite\n// it in place, preserving the paragraph's formatting.\n\nconst replacements = [\n  {\n    find: \"Play Jackpot Raiders Free | Exciting Indiana Jones-themed Slot Game\",\n    replace: \"Play Jackpot Raiders for Free\",\n  },\n  {\n    find: \"Bonus features with scatter symbols, map symbols, and Pick and Click game\",\n    replace: \"Max payout of 10,000x the stake\",\n  },\n  {\n    find: \"10,000x the stake maximum payout\",\n    replace: \"Bonus features with scatter symbols and Pick and Click game\",\n  },\n  {\n    find: \"20 paylines and x3 multiplier on all jackpot free spins\",\n    replace: \"Jackpot prizes and free spins with x3 multiplier\",\n  },\n  {\n    find: \"Only one setting for coins\",\n    replace: \"Limited coin setting\",\n  },\n  {\n    find: \"Limited range of wager amounts\",\n    replace: \"Limited wager range\",\n  },\n  {\n    find: \"Read our review of Jackpot Raiders and play this exciting Indiana Jones-style slot game for free. Collect scatter and map symbols for bonus features and a chance at winning the jackpot.\",\n    replace: \"Read our review of Jackpot Raiders and play this exciting slot game for free.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit to the Jackpot Raiders review.\n# Each change is a whole-paragraph text replacement; Find/Replace across the\n# full document body (wdReplaceAll) handles the title text that repeats\n# (once as the Heading1, once as the bold call-to-action paragraph) in a\n# single call, and is harmlessly a \"replace 1 occurrence\" for the rest.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-All \"Play Jackpot Raiders Free | Exciting Indiana Jones-themed Slot Game\" \"Play Jackpot Raiders for Free\"\n\nReplace-All \"Bonus features with scatter symbols, map symbols, and Pick and Click game\" \"Max payout of 10,000x the stake\"\n\nReplace-All \"10,000x the stake maximum payout\" \"Bonus features with scatter symbols and Pick and Click game\"\n\nReplace-All \"20 paylines and x3 multiplier on all jackpot free spins\" \"Jackpot prizes and free spins with x3 multiplier\"\n\nReplace-All \"Only one setting for coins\" \"Limited coin setting\"\n\nReplace-All \"Limited range of wager amounts\" \"Limited wager range\"\n\nReplace-All \"Read our review of Jackpot Raiders and play this exciting Indiana Jones-style slot game for free. Collect scatter and map symbols for bonus features and a chance at winning the jackpot.\" \"Read our review of Jackpot Raiders and play this exciting slot game for free.\"\n"}
